$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Q2_20_21 all data"

$new = $wb.Worksheets.Add($null, $ws)
$new.Name = "Q2_20_21 Count"

$new.Range("B3").Value = "BRD Risk Category"
$new.Range("C3").Value = "Low"
$new.Range("D3").Value = "Medium"
$new.Range("E3").Value = "High"
$new.Range("F3").Value = "Total"
$new.Range("B4").Value = "Economic"
$new.Range("C4").Value = 4
$new.Range("D4").Value = 7
$new.Range("E4").Value = 0
$new.Range("F4").Value = 11
$new.Range("B5").Value = "Technological"
$new.Range("C5").Value = 0
$new.Range("D5").Value = 3
$new.Range("E5").Value = 1
$new.Range("F5").Value = 4
$new.Range("B6").Value = "Political"
$new.Range("C6").Value = 2
$new.Range("D6").Value = 3
$new.Range("E6").Value = 0
$new.Range("F6").Value = 5
$new.Range("B7").Value = "Legal"
$new.Range("C7").Value = 1
$new.Range("D7").Value = 2
$new.Range("E7").Value = 0
$new.Range("F7").Value = 3
$new.Range("B8").Value = "Social"
$new.Range("C8").Value = 2
$new.Range("D8").Value = 1
$new.Range("E8").Value = 0
$new.Range("F8").Value = 3
$new.Range("B9").Value = "Environmental"
$new.Range("C9").Value = 0
$new.Range("D9").Value = 2
$new.Range("E9").Value = 0
$new.Range("F9").Value = 2
$new.Range("B12").Value = "BRD Primary Risk to"
$new.Range("C12").Value = "Low"
$new.Range("D12").Value = "Medium"
$new.Range("E12").Value = "High"
$new.Range("F12").Value = "Total"
$new.Range("B13").Value = "Costs"
$new.Range("C13").Value = 3
$new.Range("D13").Value = 6
$new.Range("E13").Value = 0
$new.Range("F13").Value = 9
$new.Range("B14").Value = "Schedule"
$new.Range("C14").Value = 3
$new.Range("D14").Value = 6
$new.Range("E14").Value = 0
$new.Range("F14").Value = 9
$new.Range("B15").Value = "Benefits"
$new.Range("C15").Value = 3
$new.Range("D15").Value = 6
$new.Range("E15").Value = 1
$new.Range("F15").Value = 10
$new.Range("B18").Value = "BRD Internal Control"
$new.Range("C18").Value = "Low"
$new.Range("D18").Value = "Medium"
$new.Range("E18").Value = "High"
$new.Range("F18").Value = "Total"
$new.Range("B19").Value = "Treat - Prevent"
$new.Range("C19").Value = 5
$new.Range("D19").Value = 12
$new.Range("E19").Value = 1
$new.Range("F19").Value = 18
$new.Range("B20").Value = "Treat - Directive"
$new.Range("C20").Value = 1
$new.Range("D20").Value = 0
$new.Range("E20").Value = 0
$new.Range("F20").Value = 1
$new.Range("B21").Value = "Treat - Corrective"
$new.Range("C21").Value = 2
$new.Range("D21").Value = 4
$new.Range("E21").Value = 0
$new.Range("F21").Value = 6
$new.Range("B22").Value = "Tolerate"
$new.Range("C22").Value = 1
$new.Range("D22").Value = 2
$new.Range("E22").Value = 0
$new.Range("F22").Value = 3
$new.Range("B25").Value = "BRD Residual Impact"
$new.Range("C25").Value = "Low"
$new.Range("D25").Value = "Medium"
$new.Range("E25").Value = "High"
$new.Range("F25").Value = "Total"
$new.Range("B26").Value = "High"
$new.Range("C26").Value = 0
$new.Range("D26").Value = 6
$new.Range("E26").Value = 1
$new.Range("F26").Value = 7
$new.Range("B27").Value = "Medium"
$new.Range("C27").Value = 3
$new.Range("D27").Value = 10
$new.Range("E27").Value = 0
$new.Range("F27").Value = 13
$new.Range("B28").Value = "Low"
$new.Range("C28").Value = 5
$new.Range("D28").Value = 1
$new.Range("E28").Value = 0
$new.Range("F28").Value = 6
$new.Range("B29").Value = "Very High"
$new.Range("C29").Value = 0
$new.Range("D29").Value = 1
$new.Range("E29").Value = 0
$new.Range("F29").Value = 1
$new.Range("B30").Value = "Very Low"
$new.Range("C30").Value = 1
$new.Range("D30").Value = 0
$new.Range("E30").Value = 0
$new.Range("F30").Value = 1
$new.Range("B33").Value = "BRD Residual Likelihood"
$new.Range("C33").Value = "Low"
$new.Range("D33").Value = "Medium"
$new.Range("E33").Value = "High"
$new.Range("F33").Value = "Total"
$new.Range("B34").Value = "Medium"
$new.Range("C34").Value = 1
$new.Range("D34").Value = 15
$new.Range("E34").Value = 0
$new.Range("F34").Value = 16
$new.Range("B35").Value = "Low"
$new.Range("C35").Value = 8
$new.Range("D35").Value = 2
$new.Range("E35").Value = 0
$new.Range("F35").Value = 10
$new.Range("B36").Value = "High"
$new.Range("C36").Value = 0
$new.Range("D36").Value = 1
$new.Range("E36").Value = 1
$new.Range("F36").Value = 2
$new.Range("B39").Value = "Severity Score Risk Category"
$new.Range("C39").Value = "Low"
$new.Range("D39").Value = "Medium"
$new.Range("E39").Value = "High"
$new.Range("F39").Value = "Total"
$new.Range("B40").Value = "Medium"
$new.Range("C40").Value = 0
$new.Range("D40").Value = 18
$new.Range("E40").Value = 0
$new.Range("F40").Value = 18
$new.Range("B41").Value = "Low"
$new.Range("C41").Value = 9
$new.Range("D41").Value = 0
$new.Range("E41").Value = 0
$new.Range("F41").Value = 9
$new.Range("B42").Value = "High"
$new.Range("C42").Value = 0
$new.Range("D42").Value = 0
$new.Range("E42").Value = 1
$new.Range("F42").Value = 1
$new.Range("B45").Value = "BRD Has this Risk turned into an Issue?"
$new.Range("C45").Value = "Low"
$new.Range("D45").Value = "Medium"
$new.Range("E45").Value = "High"
$new.Range("F45").Value = "Total"
$new.Range("B46").Value = "No"
$new.Range("C46").Value = 6
$new.Range("D46").Value = 15
$new.Range("E46").Value = 0
$new.Range("F46").Value = 21
$new.Range("C47").Value = 3
$new.Range("D47").Value = 3
$new.Range("E47").Value = 1
$new.Range("F47").Value = 7
